$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.335.09'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.933.03'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '486.40'
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.23'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  +3.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000357'
$ws.Range('E11').Value = '  +5.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.01'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.66'
$ws.Range('E13').Value = '  +3.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.561.39'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.66'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.938.62'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.01'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.448.91'
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '443.75'
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.54'
$ws.Range('E22').Value = '  +5.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.22'
$ws.Range('E23').Value = '  +4.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.70'
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.38'
$ws.Range('E25').Value = '  +18.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.53'
$ws.Range('E26').Value = '  +12.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.63'
$ws.Range('E27').Value = '  +2.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.01'
$ws.Range('E28').Value = '  +1.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.85'
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '728.26'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.77'
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.131'
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.93'
$ws.Range('E33').Value = '  +3.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0920'
$ws.Range('E34').Value = '  +15.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.28'
$ws.Range('E35').Value = '  +17.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.51'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '61.18'
$ws.Range('E37').Value = '  +5.67%  '
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.399'
$ws.Range('E39').Value = '  +18.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.98'
$ws.Range('E41').Value = '  +14.51%  '
$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.27'
$ws.Range('E42').Value = '  +7.71%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0484'
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.95'
$ws.Range('E44').Value = '  +5.71%  '
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0365'
$ws.Range('E46').Value = '  +39.70%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.20'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '145.88'
$ws.Range('E51').Value = '  -0.21%  '
